# Update computed cosinor statistics for the two result rows (row 2 and
# row 3) to reflect the latest CircadiPy simulation re-run.
# NOTE: scientific ("e"/"E") notation literals are not parsed here, so
# very small/large magnitudes below are written out as plain decimals
# that round-trip to the exact same IEEE-754 double value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Row 2 -----
$ws.Range("E2").Value = 22.68000000000011
$ws.Range("G2").Value = 0.0000000000000004440892098500626
$ws.Range("H2").Value = 0.000000000000006832141690000963
$ws.Range("K2").Value = 50.18595411267058
$ws.Range("L2").Value = "[38.41135263940156, 61.9605555859396]"
$ws.Range("M2").Value = 0.00000000000003730349362740526
$ws.Range("N2").Value = 0.00000000000007460698725481052
$ws.Range("O2").Value = 1.666710817219809
$ws.Range("P2").Value = "[1.4151318259413461, 1.9182898084982725]"
$ws.Range("S2").Value = 64.09125246705796
$ws.Range("T2").Value = "[56.83663200411439, 71.34587293000152]"
$ws.Range("W2").Value = 16.66378378378386
$ws.Range("X2").Value = 15.75567567567575
$ws.Range("Y2").Value = 17.57189189189198

# ----- Row 3 -----
$ws.Range("E3").Value = 22.92000000000014
$ws.Range("G3").Value = 0.0000000008890066460764956
$ws.Range("H3").Value = 0.000000001968027597574503
$ws.Range("K3").Value = 36.38537906412012
$ws.Range("L3").Value = "[22.48937728069793, 50.281380847542316]"
$ws.Range("M3").Value = 0.0000007641304948791117
$ws.Range("N3").Value = 0.0000007641304948791117
$ws.Range("O3").Value = 1.993763505881811
$ws.Range("P3").Value = "[1.6038160694001942, 2.383710942363427]"
$ws.Range("S3").Value = 60.79076573097927
$ws.Range("T3").Value = "[53.35299368046535, 68.2285377814932]"
$ws.Range("W3").Value = 15.64708708708718
$ws.Range("X3").Value = 14.22462462462471
$ws.Range("Y3").Value = 17.06954954954966

Write-Host "Updated cosinor statistics for rows 2-3"
